$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Femacal de La Calera - Haba) needs to be inserted
# as the new row 59; all existing rows 59-66 shift down to 60-67.
$ws.Rows.Item(59).Insert()

$ws.Range("A59").Value = 3
$ws.Range("B59").Value = "Femacal de La Calera"
$ws.Range("C59").Value = "Coquimbo"
$ws.Range("D59").Value = 44474
$ws.Range("E59").Value = 5
$ws.Range("F59").Value = 100112026
$ws.Range("G59").Value = "Haba"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 73
$ws.Range("K59").Value = 9000
$ws.Range("L59").Value = 9500
$ws.Range("M59").Value = 9260
$ws.Range("N59").Value = "$/malla 25 kilos"
$ws.Range("O59").Value = "Provincia de Limarí"
$ws.Range("P59").Value = 370
$ws.Range("Q59").Value = 25
$ws.Range("R59").Value = "Hortaliza"
